# VBA: 0.30 release, ready for beta testing
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Workbook now recalculates automatically instead of manually.
$excel.Calculation = -4105

# Row 38: add a "Bugs" (F) note, plus a new "Ideas" (D) note; the old
# "Tasks" (E) note on this row is replaced by the one that used to live
# further down (HC selection automation).
$ws.Cells.Item(38, 4).Value = "Availability calc properly"
$ws.Cells.Item(38, 5).Value = "HC selection automation"
$ws.Cells.Item(38, 6).Value = "noticing some performance hits, need to profile sometime"

# Row 39: add an "Ideas" (D) note.
$ws.Cells.Item(39, 4).Value = "cons per expos if needed"

# Row 40 used to only hold a "Tasks" (E) note; it now holds the next
# "Features" (B) note instead, and the old row 42 disappears.
$ws.Cells.Item(40, 2).Value = "dec sum remainings"
$ws.Cells.Item(40, 5).Clear()

# Row 41 used to only hold a "Tasks" (E) note; it now starts the 0.31
# version entry, and the old note has been relocated above.
$ws.Cells.Item(41, 1).Value = 0.31
$ws.Cells.Item(41, 5).Clear()

# Old row 42's lone "Tasks" (E) note has been relocated above; clearing
# it removes the now-empty trailing row.
$ws.Cells.Item(42, 5).Clear()
